# ============================================================
# Update betting-odds data for 'Ecuador LigaPro Serie A' sheet
# ============================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Rows that already exist: update changed cell values ----
# Row 130
$ws.Cells.Item(130, 1).Value = 128
$ws.Cells.Item(130, 2).Value = 7483081
$ws.Cells.Item(130, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(130, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(130, 5).Value = 45255.83333333334
$ws.Cells.Item(130, 6).Value = 'Deportivo Cuenca'
$ws.Cells.Item(130, 7).Value = 'El Nacional'
$ws.Cells.Item(130, 8).Value = 1
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 'H'
$ws.Cells.Item(130, 11).Value = 2.75
$ws.Cells.Item(130, 12).Value = 3.25
$ws.Cells.Item(130, 13).Value = 2.55
$ws.Cells.Item(130, 14).Value = 3
$ws.Cells.Item(130, 15).Value = 3.3
$ws.Cells.Item(130, 16).Value = 2.3
$ws.Cells.Item(130, 17).Value = 0.25
$ws.Cells.Item(130, 18).Value = 1.825
$ws.Cells.Item(130, 19).Value = 1.975
$ws.Cells.Item(130, 20).Value = 2.75
$ws.Cells.Item(130, 21).Value = 2
$ws.Cells.Item(130, 22).Value = 1.8
$ws.Cells.Item(130, 23).Value = 2
$ws.Cells.Item(130, 24).Value = -1
$ws.Cells.Item(130, 25).Value = -1
$ws.Cells.Item(130, 26).Value = 0.825
$ws.Cells.Item(130, 27).Value = -1
$ws.Cells.Item(130, 28).Value = -1
$ws.Cells.Item(130, 29).Value = 0.8

# Row 131
$ws.Cells.Item(131, 1).Value = 129
$ws.Cells.Item(131, 2).Value = 7483189
$ws.Cells.Item(131, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(131, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(131, 5).Value = 45255.83333333334
$ws.Cells.Item(131, 6).Value = 'Independiente del Valle'
$ws.Cells.Item(131, 7).Value = 'Orense'
$ws.Cells.Item(131, 8).Value = 2
$ws.Cells.Item(131, 9).Value = 2
$ws.Cells.Item(131, 10).Value = 'D'
$ws.Cells.Item(131, 11).Value = 1.4
$ws.Cells.Item(131, 12).Value = 4.75
$ws.Cells.Item(131, 13).Value = 7
$ws.Cells.Item(131, 14).Value = 1.4
$ws.Cells.Item(131, 15).Value = 4.5
$ws.Cells.Item(131, 16).Value = 8
$ws.Cells.Item(131, 17).Value = -1.25
$ws.Cells.Item(131, 18).Value = 1.875
$ws.Cells.Item(131, 19).Value = 1.925
$ws.Cells.Item(131, 20).Value = 2.5
$ws.Cells.Item(131, 21).Value = 1.925
$ws.Cells.Item(131, 22).Value = 1.875
$ws.Cells.Item(131, 23).Value = -1
$ws.Cells.Item(131, 24).Value = 3.5
$ws.Cells.Item(131, 25).Value = -1
$ws.Cells.Item(131, 26).Value = -1
$ws.Cells.Item(131, 27).Value = 0.925
$ws.Cells.Item(131, 28).Value = 0.925
$ws.Cells.Item(131, 29).Value = -1

# Row 132
$ws.Cells.Item(132, 1).Value = 130
$ws.Cells.Item(132, 2).Value = 7483281
$ws.Cells.Item(132, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(132, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(132, 5).Value = 45255.83333333334
$ws.Cells.Item(132, 6).Value = 'SD Aucas'
$ws.Cells.Item(132, 7).Value = 'Delfin SC'
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 'D'
$ws.Cells.Item(132, 11).Value = 1.909
$ws.Cells.Item(132, 12).Value = 3.25
$ws.Cells.Item(132, 13).Value = 4.2
$ws.Cells.Item(132, 14).Value = 1.909
$ws.Cells.Item(132, 15).Value = 3.5
$ws.Cells.Item(132, 16).Value = 4
$ws.Cells.Item(132, 17).Value = -0.5
$ws.Cells.Item(132, 18).Value = 1.9
$ws.Cells.Item(132, 19).Value = 1.9
$ws.Cells.Item(132, 20).Value = 2.5
$ws.Cells.Item(132, 21).Value = 1.8
$ws.Cells.Item(132, 22).Value = 2
$ws.Cells.Item(132, 23).Value = -1
$ws.Cells.Item(132, 24).Value = 2.5
$ws.Cells.Item(132, 25).Value = -1
$ws.Cells.Item(132, 26).Value = -1
$ws.Cells.Item(132, 27).Value = 0.8999999999999999
$ws.Cells.Item(132, 28).Value = -1
$ws.Cells.Item(132, 29).Value = 1

# Row 133
$ws.Cells.Item(133, 1).Value = 131
$ws.Cells.Item(133, 2).Value = 7483247
$ws.Cells.Item(133, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(133, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(133, 5).Value = 45255.83333333334
$ws.Cells.Item(133, 6).Value = 'Mushuc Runa'
$ws.Cells.Item(133, 7).Value = 'Universidad Catolica del Ecuador'
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 2
$ws.Cells.Item(133, 10).Value = 'A'
$ws.Cells.Item(133, 11).Value = 3.25
$ws.Cells.Item(133, 12).Value = 3.2
$ws.Cells.Item(133, 13).Value = 2.25
$ws.Cells.Item(133, 14).Value = 3.5
$ws.Cells.Item(133, 15).Value = 3.25
$ws.Cells.Item(133, 16).Value = 2.1
$ws.Cells.Item(133, 17).Value = 0.5
$ws.Cells.Item(133, 18).Value = 1.775
$ws.Cells.Item(133, 19).Value = 2.025
$ws.Cells.Item(133, 20).Value = 2.5
$ws.Cells.Item(133, 21).Value = 1.9
$ws.Cells.Item(133, 22).Value = 1.9
$ws.Cells.Item(133, 23).Value = -1
$ws.Cells.Item(133, 24).Value = -1
$ws.Cells.Item(133, 25).Value = 1.1
$ws.Cells.Item(133, 26).Value = -1
$ws.Cells.Item(133, 27).Value = 1.025
$ws.Cells.Item(133, 28).Value = -1
$ws.Cells.Item(133, 29).Value = 0.8999999999999999

# Row 134
$ws.Cells.Item(134, 1).Value = 132
$ws.Cells.Item(134, 2).Value = 7482867
$ws.Cells.Item(134, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(134, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(134, 5).Value = 45256.83333333334
$ws.Cells.Item(134, 6).Value = 'Cumbaya FC'
$ws.Cells.Item(134, 7).Value = 'LDU Quito'
$ws.Cells.Item(134, 8).Value = 1
$ws.Cells.Item(134, 9).Value = 2
$ws.Cells.Item(134, 10).Value = 'A'
$ws.Cells.Item(134, 11).Value = 5.25
$ws.Cells.Item(134, 12).Value = 3.75
$ws.Cells.Item(134, 13).Value = 1.65
$ws.Cells.Item(134, 14).Value = 9
$ws.Cells.Item(134, 15).Value = 4.5
$ws.Cells.Item(134, 16).Value = 1.363
$ws.Cells.Item(134, 17).Value = 1.25
$ws.Cells.Item(134, 18).Value = 1.975
$ws.Cells.Item(134, 19).Value = 1.825
$ws.Cells.Item(134, 20).Value = 2.5
$ws.Cells.Item(134, 21).Value = 1.825
$ws.Cells.Item(134, 22).Value = 1.975
$ws.Cells.Item(134, 23).Value = -1
$ws.Cells.Item(134, 24).Value = -1
$ws.Cells.Item(134, 25).Value = 0.363
$ws.Cells.Item(134, 26).Value = 0.4875
$ws.Cells.Item(134, 27).Value = -0.5
$ws.Cells.Item(134, 28).Value = 0.825
$ws.Cells.Item(134, 29).Value = -1

# Row 135
$ws.Cells.Item(135, 1).Value = 133
$ws.Cells.Item(135, 2).Value = 7483188
$ws.Cells.Item(135, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(135, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(135, 5).Value = 45256.83333333334
$ws.Cells.Item(135, 6).Value = 'Gualaceo SC'
$ws.Cells.Item(135, 7).Value = 'Emelec'
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 9).Value = 2
$ws.Cells.Item(135, 10).Value = 'A'
$ws.Cells.Item(135, 11).Value = 3.6
$ws.Cells.Item(135, 12).Value = 3.3
$ws.Cells.Item(135, 13).Value = 2.05
$ws.Cells.Item(135, 14).Value = 2.6
$ws.Cells.Item(135, 15).Value = 3.25
$ws.Cells.Item(135, 16).Value = 2.75
$ws.Cells.Item(135, 17).Value = 0
$ws.Cells.Item(135, 18).Value = 1.8
$ws.Cells.Item(135, 19).Value = 2
$ws.Cells.Item(135, 20).Value = 2.5
$ws.Cells.Item(135, 21).Value = 1.975
$ws.Cells.Item(135, 22).Value = 1.825
$ws.Cells.Item(135, 23).Value = -1
$ws.Cells.Item(135, 24).Value = -1
$ws.Cells.Item(135, 25).Value = 1.75
$ws.Cells.Item(135, 26).Value = -1
$ws.Cells.Item(135, 27).Value = 1
$ws.Cells.Item(135, 28).Value = -1
$ws.Cells.Item(135, 29).Value = 0.825

# Row 136
$ws.Cells.Item(136, 1).Value = 134
$ws.Cells.Item(136, 2).Value = 7483306
$ws.Cells.Item(136, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(136, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(136, 5).Value = 45256.83333333334
$ws.Cells.Item(136, 6).Value = 'Tecnico Universitario'
$ws.Cells.Item(136, 7).Value = 'Club Atletico Libertad'
$ws.Cells.Item(136, 8).Value = 1
$ws.Cells.Item(136, 9).Value = 1
$ws.Cells.Item(136, 10).Value = 'D'
$ws.Cells.Item(136, 11).Value = 1.5
$ws.Cells.Item(136, 12).Value = 4.333
$ws.Cells.Item(136, 13).Value = 5.75
$ws.Cells.Item(136, 14).Value = 1.533
$ws.Cells.Item(136, 15).Value = 4.2
$ws.Cells.Item(136, 16).Value = 5.5
$ws.Cells.Item(136, 17).Value = -1
$ws.Cells.Item(136, 18).Value = 1.925
$ws.Cells.Item(136, 19).Value = 1.875
$ws.Cells.Item(136, 20).Value = 2.25
$ws.Cells.Item(136, 21).Value = 1.8
$ws.Cells.Item(136, 22).Value = 2
$ws.Cells.Item(136, 23).Value = -1
$ws.Cells.Item(136, 24).Value = 3.2
$ws.Cells.Item(136, 25).Value = -1
$ws.Cells.Item(136, 26).Value = -1
$ws.Cells.Item(136, 27).Value = 0.875
$ws.Cells.Item(136, 28).Value = -0.5
$ws.Cells.Item(136, 29).Value = 0.5

# Row 137
$ws.Cells.Item(137, 1).Value = 135
$ws.Cells.Item(137, 2).Value = 7482832
$ws.Cells.Item(137, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(137, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(137, 5).Value = 45256.83333333334
$ws.Cells.Item(137, 6).Value = 'Barcelona Guayaquil'
$ws.Cells.Item(137, 7).Value = 'Guayaquil City'
$ws.Cells.Item(137, 8).Value = 2
$ws.Cells.Item(137, 9).Value = 1
$ws.Cells.Item(137, 10).Value = 'H'
$ws.Cells.Item(137, 11).Value = 1.363
$ws.Cells.Item(137, 12).Value = 5
$ws.Cells.Item(137, 13).Value = 7.5
$ws.Cells.Item(137, 14).Value = 1.444
$ws.Cells.Item(137, 15).Value = 4
$ws.Cells.Item(137, 16).Value = 8
$ws.Cells.Item(137, 17).Value = -1.25
$ws.Cells.Item(137, 18).Value = 2.05
$ws.Cells.Item(137, 19).Value = 1.75
$ws.Cells.Item(137, 20).Value = 2.5
$ws.Cells.Item(137, 21).Value = 1.95
$ws.Cells.Item(137, 22).Value = 1.85
$ws.Cells.Item(137, 23).Value = 0.444
$ws.Cells.Item(137, 24).Value = -1
$ws.Cells.Item(137, 25).Value = -1
$ws.Cells.Item(137, 26).Value = -0.5
$ws.Cells.Item(137, 27).Value = 0.375
$ws.Cells.Item(137, 28).Value = 0.95
$ws.Cells.Item(137, 29).Value = -1

# Row 142
$ws.Cells.Item(142, 1).Value = 140
$ws.Cells.Item(142, 2).Value = 7528858
$ws.Cells.Item(142, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(142, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(142, 5).Value = 45263.83333333334
$ws.Cells.Item(142, 6).Value = 'Orense'
$ws.Cells.Item(142, 7).Value = 'SD Aucas'
$ws.Cells.Item(142, 8).Value = 1
$ws.Cells.Item(142, 9).Value = 2
$ws.Cells.Item(142, 10).Value = 'A'
$ws.Cells.Item(142, 11).Value = 2.2
$ws.Cells.Item(142, 12).Value = 3.2
$ws.Cells.Item(142, 13).Value = 3.2
$ws.Cells.Item(142, 14).Value = 1.95
$ws.Cells.Item(142, 15).Value = 3.2
$ws.Cells.Item(142, 16).Value = 3.8
$ws.Cells.Item(142, 17).Value = -0.5
$ws.Cells.Item(142, 18).Value = 1.95
$ws.Cells.Item(142, 19).Value = 1.85
$ws.Cells.Item(142, 20).Value = 2.25
$ws.Cells.Item(142, 21).Value = 1.85
$ws.Cells.Item(142, 22).Value = 1.95
$ws.Cells.Item(142, 23).Value = -1
$ws.Cells.Item(142, 24).Value = -1
$ws.Cells.Item(142, 25).Value = 2.8
$ws.Cells.Item(142, 26).Value = -1
$ws.Cells.Item(142, 27).Value = 0.8500000000000001
$ws.Cells.Item(142, 28).Value = 0.8500000000000001
$ws.Cells.Item(142, 29).Value = -1

# Row 143
$ws.Cells.Item(143, 1).Value = 141
$ws.Cells.Item(143, 2).Value = 7528852
$ws.Cells.Item(143, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(143, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(143, 5).Value = 45263.83333333334
$ws.Cells.Item(143, 6).Value = 'Delfin SC'
$ws.Cells.Item(143, 7).Value = 'Tecnico Universitario'
$ws.Cells.Item(143, 8).Value = 2
$ws.Cells.Item(143, 9).Value = 2
$ws.Cells.Item(143, 10).Value = 'D'
$ws.Cells.Item(143, 11).Value = 2.1
$ws.Cells.Item(143, 12).Value = 3.4
$ws.Cells.Item(143, 13).Value = 3.1
$ws.Cells.Item(143, 14).Value = 2.1
$ws.Cells.Item(143, 15).Value = 3.4
$ws.Cells.Item(143, 16).Value = 3.1
$ws.Cells.Item(143, 17).Value = -0.25
$ws.Cells.Item(143, 18).Value = 1.8
$ws.Cells.Item(143, 19).Value = 2
$ws.Cells.Item(143, 20).Value = 2.25
$ws.Cells.Item(143, 21).Value = 1.9
$ws.Cells.Item(143, 22).Value = 1.9
$ws.Cells.Item(143, 23).Value = -1
$ws.Cells.Item(143, 24).Value = 2.4
$ws.Cells.Item(143, 25).Value = -1
$ws.Cells.Item(143, 26).Value = -0.5
$ws.Cells.Item(143, 27).Value = 0.5
$ws.Cells.Item(143, 28).Value = 0.8999999999999999
$ws.Cells.Item(143, 29).Value = -1

# Row 145
$ws.Cells.Item(145, 1).Value = 143
$ws.Cells.Item(145, 2).Value = 7528848
$ws.Cells.Item(145, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(145, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(145, 5).Value = 45263.83333333334
$ws.Cells.Item(145, 6).Value = 'Emelec'
$ws.Cells.Item(145, 7).Value = 'Deportivo Cuenca'
$ws.Cells.Item(145, 8).Value = 2
$ws.Cells.Item(145, 9).Value = 1
$ws.Cells.Item(145, 10).Value = 'H'
$ws.Cells.Item(145, 11).Value = 1.75
$ws.Cells.Item(145, 12).Value = 3.5
$ws.Cells.Item(145, 13).Value = 4.2
$ws.Cells.Item(145, 14).Value = 2.4
$ws.Cells.Item(145, 15).Value = 3.1
$ws.Cells.Item(145, 16).Value = 2.75
$ws.Cells.Item(145, 17).Value = -0.25
$ws.Cells.Item(145, 18).Value = 2.05
$ws.Cells.Item(145, 19).Value = 1.75
$ws.Cells.Item(145, 20).Value = 2.25
$ws.Cells.Item(145, 21).Value = 1.8
$ws.Cells.Item(145, 22).Value = 2
$ws.Cells.Item(145, 23).Value = 1.4
$ws.Cells.Item(145, 24).Value = -1
$ws.Cells.Item(145, 25).Value = -1
$ws.Cells.Item(145, 26).Value = 1.05
$ws.Cells.Item(145, 27).Value = -1
$ws.Cells.Item(145, 28).Value = 0.8
$ws.Cells.Item(145, 29).Value = -1

# Row 200
$ws.Cells.Item(200, 1).Value = 198
$ws.Cells.Item(200, 2).Value = 7773501
$ws.Cells.Item(200, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(200, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(200, 5).Value = 45394.875
$ws.Cells.Item(200, 6).Value = 'Deportivo Cuenca'
$ws.Cells.Item(200, 7).Value = 'SD Aucas'
$ws.Cells.Item(200, 8).Value = 3
$ws.Cells.Item(200, 9).Value = 3
$ws.Cells.Item(200, 10).Value = 'D'
$ws.Cells.Item(200, 11).Value = 3
$ws.Cells.Item(200, 12).Value = 3.2
$ws.Cells.Item(200, 13).Value = 2.25
$ws.Cells.Item(200, 14).Value = 3
$ws.Cells.Item(200, 15).Value = 3.1
$ws.Cells.Item(200, 16).Value = 2.25
$ws.Cells.Item(200, 17).Value = 0.25
$ws.Cells.Item(200, 18).Value = 1.8
$ws.Cells.Item(200, 19).Value = 2
$ws.Cells.Item(200, 20).Value = 2.25
$ws.Cells.Item(200, 21).Value = 1.95
$ws.Cells.Item(200, 22).Value = 1.85
$ws.Cells.Item(200, 23).Value = -1
$ws.Cells.Item(200, 24).Value = 2.1
$ws.Cells.Item(200, 25).Value = -1
$ws.Cells.Item(200, 26).Value = 0.4
$ws.Cells.Item(200, 27).Value = -0.5
$ws.Cells.Item(200, 28).Value = 0.95
$ws.Cells.Item(200, 29).Value = -1

# Row 201
$ws.Cells.Item(201, 1).Value = 199
$ws.Cells.Item(201, 2).Value = 7773503
$ws.Cells.Item(201, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(201, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(201, 5).Value = 45395.625
$ws.Cells.Item(201, 6).Value = 'Mushuc Runa'
$ws.Cells.Item(201, 7).Value = 'Tecnico Universitario'
$ws.Cells.Item(201, 8).Value = 1
$ws.Cells.Item(201, 9).Value = 2
$ws.Cells.Item(201, 10).Value = 'A'
$ws.Cells.Item(201, 11).Value = 2.1
$ws.Cells.Item(201, 12).Value = 3.1
$ws.Cells.Item(201, 13).Value = 3.4
$ws.Cells.Item(201, 14).Value = 2.25
$ws.Cells.Item(201, 15).Value = 3
$ws.Cells.Item(201, 16).Value = 3.1
$ws.Cells.Item(201, 17).Value = -0.25
$ws.Cells.Item(201, 18).Value = 1.95
$ws.Cells.Item(201, 19).Value = 1.85
$ws.Cells.Item(201, 20).Value = 2.25
$ws.Cells.Item(201, 21).Value = 1.925
$ws.Cells.Item(201, 22).Value = 1.875
$ws.Cells.Item(201, 23).Value = -1
$ws.Cells.Item(201, 24).Value = -1
$ws.Cells.Item(201, 25).Value = 2.1
$ws.Cells.Item(201, 26).Value = -1
$ws.Cells.Item(201, 27).Value = 0.8500000000000001
$ws.Cells.Item(201, 28).Value = 0.925
$ws.Cells.Item(201, 29).Value = -1

# Row 202
$ws.Cells.Item(202, 1).Value = 200
$ws.Cells.Item(202, 2).Value = 8069719
$ws.Cells.Item(202, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(202, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(202, 5).Value = 45395.72916666666
$ws.Cells.Item(202, 6).Value = 'Macara'
$ws.Cells.Item(202, 7).Value = 'Orense'
$ws.Cells.Item(202, 8).Value = 0
$ws.Cells.Item(202, 9).Value = 0
$ws.Cells.Item(202, 10).Value = 'D'
$ws.Cells.Item(202, 11).Value = 1.95
$ws.Cells.Item(202, 12).Value = 3.25
$ws.Cells.Item(202, 13).Value = 3.5
$ws.Cells.Item(202, 14).Value = 1.571
$ws.Cells.Item(202, 15).Value = 3.6
$ws.Cells.Item(202, 16).Value = 5.25
$ws.Cells.Item(202, 17).Value = -1
$ws.Cells.Item(202, 18).Value = 2
$ws.Cells.Item(202, 19).Value = 1.8
$ws.Cells.Item(202, 20).Value = 2.5
$ws.Cells.Item(202, 21).Value = 2
$ws.Cells.Item(202, 22).Value = 1.8
$ws.Cells.Item(202, 23).Value = -1
$ws.Cells.Item(202, 24).Value = 2.6
$ws.Cells.Item(202, 25).Value = -1
$ws.Cells.Item(202, 26).Value = -1
$ws.Cells.Item(202, 27).Value = 0.8
$ws.Cells.Item(202, 28).Value = -1
$ws.Cells.Item(202, 29).Value = 0.8

# ---- New rows appended at the end (203-208) ----
# Copy formatting for column A (row label style) and column E (date style)
# from an existing data row so the new rows look consistent.
$ws.Range('A130').Copy()
$ws.Range('A203').PasteSpecial(-4122)
$ws.Range('E130').Copy()
$ws.Range('E203').PasteSpecial(-4122)
$ws.Range('A130').Copy()
$ws.Range('A204').PasteSpecial(-4122)
$ws.Range('E130').Copy()
$ws.Range('E204').PasteSpecial(-4122)
$ws.Range('A130').Copy()
$ws.Range('A205').PasteSpecial(-4122)
$ws.Range('E130').Copy()
$ws.Range('E205').PasteSpecial(-4122)
$ws.Range('A130').Copy()
$ws.Range('A206').PasteSpecial(-4122)
$ws.Range('E130').Copy()
$ws.Range('E206').PasteSpecial(-4122)
$ws.Range('A130').Copy()
$ws.Range('A207').PasteSpecial(-4122)
$ws.Range('E130').Copy()
$ws.Range('E207').PasteSpecial(-4122)
$ws.Range('A130').Copy()
$ws.Range('A208').PasteSpecial(-4122)
$ws.Range('E130').Copy()
$ws.Range('E208').PasteSpecial(-4122)

# Row 203
$ws.Cells.Item(203, 1).Value = 201
$ws.Cells.Item(203, 2).Value = 8069537
$ws.Cells.Item(203, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(203, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(203, 5).Value = 45395.83333333334
$ws.Cells.Item(203, 6).Value = 'Emelec'
$ws.Cells.Item(203, 7).Value = 'Cumbaya FC'
$ws.Cells.Item(203, 8).Value = 2
$ws.Cells.Item(203, 9).Value = 0
$ws.Cells.Item(203, 10).Value = 'H'
$ws.Cells.Item(203, 11).Value = 1.28
$ws.Cells.Item(203, 12).Value = 5.5
$ws.Cells.Item(203, 13).Value = 8.5
$ws.Cells.Item(203, 14).Value = 1.25
$ws.Cells.Item(203, 15).Value = 5.5
$ws.Cells.Item(203, 16).Value = 10
$ws.Cells.Item(203, 17).Value = -1.5
$ws.Cells.Item(203, 18).Value = 1.825
$ws.Cells.Item(203, 19).Value = 1.975
$ws.Cells.Item(203, 20).Value = 2.75
$ws.Cells.Item(203, 21).Value = 1.95
$ws.Cells.Item(203, 22).Value = 1.85
$ws.Cells.Item(203, 23).Value = 0.25
$ws.Cells.Item(203, 24).Value = -1
$ws.Cells.Item(203, 25).Value = -1
$ws.Cells.Item(203, 26).Value = 0.825
$ws.Cells.Item(203, 27).Value = -1
$ws.Cells.Item(203, 28).Value = -1
$ws.Cells.Item(203, 29).Value = 0.8500000000000001

# Row 204
$ws.Cells.Item(204, 1).Value = 202
$ws.Cells.Item(204, 2).Value = 7773511
$ws.Cells.Item(204, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(204, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(204, 5).Value = 45399.77083333334
$ws.Cells.Item(204, 6).Value = 'Cumbaya FC'
$ws.Cells.Item(204, 7).Value = 'Club Atletico Libertad'
$ws.Cells.Item(204, 11).Value = 2.4
$ws.Cells.Item(204, 12).Value = 3
$ws.Cells.Item(204, 13).Value = 3.2
$ws.Cells.Item(204, 14).Value = 2.4
$ws.Cells.Item(204, 15).Value = 3
$ws.Cells.Item(204, 16).Value = 3.2
$ws.Cells.Item(204, 17).Value = -0.25
$ws.Cells.Item(204, 18).Value = 2.025
$ws.Cells.Item(204, 19).Value = 1.775
$ws.Cells.Item(204, 20).Value = 2.25
$ws.Cells.Item(204, 21).Value = 1.9
$ws.Cells.Item(204, 22).Value = 1.9
$ws.Cells.Item(204, 23).Value = 0
$ws.Cells.Item(204, 24).Value = 0
$ws.Cells.Item(204, 25).Value = 0
$ws.Cells.Item(204, 26).Value = 0
$ws.Cells.Item(204, 27).Value = 0

# Row 205
$ws.Cells.Item(205, 1).Value = 203
$ws.Cells.Item(205, 2).Value = 7773510
$ws.Cells.Item(205, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(205, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(205, 5).Value = 45399.875
$ws.Cells.Item(205, 6).Value = 'Tecnico Universitario'
$ws.Cells.Item(205, 7).Value = 'Deportivo Cuenca'
$ws.Cells.Item(205, 11).Value = 1.95
$ws.Cells.Item(205, 12).Value = 3.25
$ws.Cells.Item(205, 13).Value = 4.2
$ws.Cells.Item(205, 14).Value = 1.95
$ws.Cells.Item(205, 15).Value = 3.25
$ws.Cells.Item(205, 16).Value = 4.2
$ws.Cells.Item(205, 17).Value = -0.5
$ws.Cells.Item(205, 18).Value = 1.95
$ws.Cells.Item(205, 19).Value = 1.85
$ws.Cells.Item(205, 20).Value = 2.25
$ws.Cells.Item(205, 21).Value = 1.9
$ws.Cells.Item(205, 22).Value = 1.9
$ws.Cells.Item(205, 23).Value = 0
$ws.Cells.Item(205, 24).Value = 0
$ws.Cells.Item(205, 25).Value = 0
$ws.Cells.Item(205, 26).Value = 0
$ws.Cells.Item(205, 27).Value = 0

# Row 206
$ws.Cells.Item(206, 1).Value = 204
$ws.Cells.Item(206, 2).Value = 7773068
$ws.Cells.Item(206, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(206, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(206, 5).Value = 45400.77083333334
$ws.Cells.Item(206, 6).Value = 'Independiente del Valle'
$ws.Cells.Item(206, 7).Value = 'Mushuc Runa'
$ws.Cells.Item(206, 11).Value = 1.333
$ws.Cells.Item(206, 12).Value = 5.5
$ws.Cells.Item(206, 13).Value = 8
$ws.Cells.Item(206, 14).Value = 1.333
$ws.Cells.Item(206, 15).Value = 5.5
$ws.Cells.Item(206, 16).Value = 8
$ws.Cells.Item(206, 17).Value = -1.5
$ws.Cells.Item(206, 18).Value = 1.925
$ws.Cells.Item(206, 19).Value = 1.875
$ws.Cells.Item(206, 20).Value = 2.75
$ws.Cells.Item(206, 21).Value = 1.8
$ws.Cells.Item(206, 22).Value = 2
$ws.Cells.Item(206, 23).Value = 0
$ws.Cells.Item(206, 24).Value = 0
$ws.Cells.Item(206, 25).Value = 0
$ws.Cells.Item(206, 26).Value = 0
$ws.Cells.Item(206, 27).Value = 0

# Row 207
$ws.Cells.Item(207, 1).Value = 205
$ws.Cells.Item(207, 2).Value = 7773508
$ws.Cells.Item(207, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(207, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(207, 5).Value = 45400.875
$ws.Cells.Item(207, 6).Value = 'Barcelona Guayaquil'
$ws.Cells.Item(207, 7).Value = 'El Nacional'
$ws.Cells.Item(207, 11).Value = 1.45
$ws.Cells.Item(207, 12).Value = 4.5
$ws.Cells.Item(207, 13).Value = 6.5
$ws.Cells.Item(207, 14).Value = 1.45
$ws.Cells.Item(207, 15).Value = 4.5
$ws.Cells.Item(207, 16).Value = 6.5
$ws.Cells.Item(207, 17).Value = -1.25
$ws.Cells.Item(207, 18).Value = 1.975
$ws.Cells.Item(207, 19).Value = 1.825
$ws.Cells.Item(207, 20).Value = 2.75
$ws.Cells.Item(207, 21).Value = 1.8
$ws.Cells.Item(207, 22).Value = 2
$ws.Cells.Item(207, 23).Value = 0
$ws.Cells.Item(207, 24).Value = 0
$ws.Cells.Item(207, 25).Value = 0
$ws.Cells.Item(207, 26).Value = 0
$ws.Cells.Item(207, 27).Value = 0

# Row 208
$ws.Cells.Item(208, 1).Value = 206
$ws.Cells.Item(208, 2).Value = 7773507
$ws.Cells.Item(208, 3).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(208, 4).Value = 'Ecuador LigaPro Serie A'
$ws.Cells.Item(208, 5).Value = 45401.77083333334
$ws.Cells.Item(208, 6).Value = 'Orense'
$ws.Cells.Item(208, 7).Value = 'LDU Quito'
$ws.Cells.Item(208, 11).Value = 4
$ws.Cells.Item(208, 12).Value = 3.25
$ws.Cells.Item(208, 13).Value = 1.909
$ws.Cells.Item(208, 14).Value = 3.8
$ws.Cells.Item(208, 15).Value = 3.25
$ws.Cells.Item(208, 16).Value = 1.95
$ws.Cells.Item(208, 17).Value = 0.5
$ws.Cells.Item(208, 18).Value = 1.825
$ws.Cells.Item(208, 19).Value = 1.975
$ws.Cells.Item(208, 20).Value = 2.75
$ws.Cells.Item(208, 21).Value = 1.975
$ws.Cells.Item(208, 22).Value = 1.825
$ws.Cells.Item(208, 23).Value = 0
$ws.Cells.Item(208, 24).Value = 0
$ws.Cells.Item(208, 25).Value = 0
$ws.Cells.Item(208, 26).Value = 0
$ws.Cells.Item(208, 27).Value = 0

$excel.CutCopyMode = 0

